$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price updates that remain plain text (values containing two or more
# dots, e.g. thousand-separated style, are never number-like so a direct
# Value assignment keeps them as text without touching any cell formatting).
$dTextValues = @{
    2 = "63.268.82"
    3 = "3.059.22"
    8 = "3.055.42"
    15 = "3.563.43"
    16 = "63.314.52"
    17 = "3.065.74"
    40 = "3.047.07"
    50 = "0.0₃0509"
}
foreach ($row in $dTextValues.Keys) {
    $ws.Range("D$row").Value = $dTextValues[$row]
}

# D-column price updates that look like plain numbers (e.g. "549.77"). Excel
# auto-converts a numeric-looking string assigned to .Value into a real number,
# which would lose the exact decimal text from the source data (and the
# trailing zeros, e.g. "1.00" -> 1). Force the cell to Text format first so the
# literal digits survive untouched, exactly like the published sheet.
$dNumericTextValues = @{
    5 = "549.77"
    6 = "140.41"
    7 = "1.00"
    9 = "0.501"
    10 = "6.50"
    11 = "0.152"
    12 = "0.455"
    14 = "34.74"
    19 = "6.76"
    20 = "483.11"
    21 = "13.81"
    22 = "0.675"
    23 = "7.25"
    24 = "81.03"
    25 = "12.71"
    26 = "0.999"
    27 = "2.76"
    28 = "7.85"
    31 = "26.14"
    33 = "2.44"
    34 = "5.66"
    35 = "55.36"
    36 = "5.98"
    37 = "462.28"
    38 = "0.0822"
    39 = "0.0397"
    41 = "0.119"
    42 = "8.21"
    43 = "2.56"
    44 = "27.77"
    45 = "0.255"
    47 = "2.04"
    49 = "116.94"
    51 = "2.07"
}
foreach ($row in $dNumericTextValues.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dNumericTextValues[$row]
}

# E-column (hourly volume/change) updates. These are always padded percentage
# strings (e.g. "  +3.20%  "), never numeric, so a direct Value assignment is safe.
$eValues = @{
    2 = "  +3.20%  "
    3 = "  +2.27%  "
    4 = "  +0.08%  "
    5 = "  +2.42%  "
    6 = "  +4.51%  "
    7 = "  -0.02%  "
    8 = "  +2.33%  "
    9 = "  +1.33%  "
    10 = "  +6.06%  "
    11 = "  +2.84%  "
    12 = "  +2.06%  "
    13 = "  +2.80%  "
    14 = "  +2.47%  "
    15 = "  +2.51%  "
    16 = "  +3.30%  "
    17 = "  +2.40%  "
    18 = "  -0.97%  "
    19 = "  +2.22%  "
    20 = "  +3.83%  "
    21 = "  +4.68%  "
    22 = "  +0.09%  "
    23 = "  +4.65%  "
    24 = "  +1.22%  "
    25 = "  +6.33%  "
    26 = "  -0.12%  "
    27 = "  +3.33%  "
    28 = "  +1.14%  "
    29 = "  +6.83%  "
    30 = "  +0.11%  "
    31 = "  +2.30%  "
    32 = "  +0.91%  "
    33 = "  +7.50%  "
    34 = "  +3.41%  "
    35 = "  +1.25%  "
    36 = "  +1.56%  "
    37 = "  +2.55%  "
    38 = "  +4.52%  "
    39 = "  +3.37%  "
    40 = "  -3.57%  "
    41 = "  -0.84%  "
    42 = "  +1.14%  "
    43 = "  +4.37%  "
    44 = "  +3.19%  "
    45 = "  +4.44%  "
    47 = "  +2.29%  "
    48 = "  +2.35%  "
    49 = "  -1.81%  "
    50 = "  +3.01%  "
    51 = "  +3.89%  "
}
foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = $eValues[$row]
}
